$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-31 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-01 Saturday", 2) | Out-Null
$d.Content.Find.Execute("693÷3=231, 0", $true, $false, $false, $false, $false, $true, 1, $false, "734÷5=146, 4", 2) | Out-Null
$d.Content.Find.Execute("964÷7=137, 5", $true, $false, $false, $false, $false, $true, 1, $false, "352÷3=117, 1", 2) | Out-Null
$d.Content.Find.Execute("118÷9=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "275÷4=68, 3", 2) | Out-Null
$d.Content.Find.Execute("581÷5=116, 1", $true, $false, $false, $false, $false, $true, 1, $false, "946÷8=118, 2", 2) | Out-Null
$d.Content.Find.Execute("912÷2=456, 0", $true, $false, $false, $false, $false, $true, 1, $false, "216÷4=54, 0", 2) | Out-Null
$d.Content.Find.Execute("258÷4=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "917÷6=152, 5", 2) | Out-Null
$d.Content.Find.Execute("682÷5=136, 2", $true, $false, $false, $false, $false, $true, 1, $false, "912÷5=182, 2", 2) | Out-Null
$d.Content.Find.Execute("382÷3=127, 1", $true, $false, $false, $false, $false, $true, 1, $false, "719÷4=179, 3", 2) | Out-Null
$d.Content.Find.Execute("850÷7=121, 3", $true, $false, $false, $false, $false, $true, 1, $false, "659÷5=131, 4", 2) | Out-Null
$d.Content.Find.Execute("808÷6=134, 4", $true, $false, $false, $false, $false, $true, 1, $false, "366÷2=183, 0", 2) | Out-Null
$d.Content.Find.Execute("394÷7=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "350÷2=175, 0", 2) | Out-Null
$d.Content.Find.Execute("807÷3=269, 0", $true, $false, $false, $false, $false, $true, 1, $false, "857÷8=107, 1", 2) | Out-Null
$d.Content.Find.Execute("671÷3=223, 2", $true, $false, $false, $false, $false, $true, 1, $false, "247÷3=82, 1", 2) | Out-Null
$d.Content.Find.Execute("720÷2=360, 0", $true, $false, $false, $false, $false, $true, 1, $false, "909÷4=227, 1", 2) | Out-Null
$d.Content.Find.Execute("487÷9=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "443÷9=49, 2", 2) | Out-Null
$d.Content.Find.Execute("853÷6=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "250÷4=62, 2", 2) | Out-Null
$d.Content.Find.Execute("545÷7=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "437÷7=62, 3", 2) | Out-Null
$d.Content.Find.Execute("933÷2=466, 1", $true, $false, $false, $false, $false, $true, 1, $false, "106÷4=26, 2", 2) | Out-Null
$d.Content.Find.Execute("337÷2=168, 1", $true, $false, $false, $false, $false, $true, 1, $false, "376÷3=125, 1", 2) | Out-Null
$d.Content.Find.Execute("991÷5=198, 1", $true, $false, $false, $false, $false, $true, 1, $false, "797÷9=88, 5", 2) | Out-Null
$d.Content.Find.Execute("408÷5=81, 3", $true, $false, $false, $false, $false, $true, 1, $false, "430÷8=53, 6", 2) | Out-Null
$d.Content.Find.Execute("647÷8=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "266÷9=29, 5", 2) | Out-Null
$d.Content.Find.Execute("116÷2=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "650÷4=162, 2", 2) | Out-Null
$d.Content.Find.Execute("925÷3=308, 1", $true, $false, $false, $false, $false, $true, 1, $false, "495÷6=82, 3", 2) | Out-Null
$d.Content.Find.Execute("674÷9=74, 8", $true, $false, $false, $false, $false, $true, 1, $false, "678÷2=339, 0", 2) | Out-Null
